# Updated the test log with more pictures
# Appends 18 new testing-log rows (21-38) to Sheet1 of the eye-detection
# testing log, each with a picture URL (as a hyperlink in column A),
# success/error status, description, the shared test date, and the
# updater's name - mirroring the existing rows 2-20 layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# (row, URL, Success/Error, Description, DateSerial, Updated by)
$rows = @(
    @(21, "http://faceresearch.org/uploads/base/african_male.jpg", "Success", "None", 41397, "Brian Nguyen"),
    @(22, "http://www.co.bibb.ga.us/TaxAssessors/images/JGordonMask.jpg", "Error", "Part of eyebrow and mustache elected", 41397, "Brian Nguyen"),
    @(23, "http://www.standardmedia.co.ke/images/friday/Angola250113.jpg", "Success", "None", 41397, "Brian Nguyen"),
    @(24, "http://everyafricanwoman.files.wordpress.com/2012/08/575034_429235650450940_156056494_n.jpg", "Error", "Eye detection window does not appear", 41397, "Brian Nguyen"),
    @(25, "http://maricopa360.com/wp-content/uploads/2009/03/raymon-tapia.jpg", "Success", "None", 41397, "Brian Nguyen"),
    @(26, "http://us.123rf.com/400wm/400/400/barsik/barsik0707/barsik070700006/1200772-portrait-of-a-young-mexican-girl-laughing.jpg", "Error", "Detects corner of mouth", 41397, "Brian Nguyen"),
    @(27, "http://www.goworldtravel.com/june05/j0227709.jpg", "Error", "Only detects right eye", 41397, "Brian Nguyen"),
    @(28, "http://1.bp.blogspot.com/-L2PBaae9jOs/TdKFi09IaSI/AAAAAAAAAKE/1Y3ymlfNDeE/s1600/indian+male+4+whysoindian.jpg", "Success", "None", 41397, "Brian Nguyen"),
    @(29, "http://faceresearch.org/uploads/base/eastasian_male.jpg", "Success", "None", 41397, "Brian Nguyen"),
    @(30, "http://0.tqn.com/d/menshair/1/0/H/A/-/-/01.jpg", "Success", "None", 41397, "Brian Nguyen"),
    @(31, "http://asianfilmfestla.org/2012/wp-content/uploads/2012/04/2012_lum_debbie_seeking_asian_female.jpg", "Success", "None", 41397, "Brian Nguyen"),
    @(32, "http://www.beautyanalysis.com/images/PG-45B---RF---Asian-Female-.jpg", "Success", "None", 41397, "Brian Nguyen"),
    @(33, "http://noahsdad.com/wp-content/2012/05/baby-down-syndrome-face-boy-smiling-9.jpg", "Error", "Detects part of the forehead", 41397, "Brian Nguyen"),
    @(34, "http://anthro.palomar.edu/abnormal/images/Down_Syndrome_child.jpg", "Error", "Doesn't detect eyes", 41397, "Brian Nguyen"),
    @(35, "http://www.pyroenergen.com/articles07/images/d-syndrome2.jpg", "Success", "None", 41397, "Brian Nguyen"),
    @(36, "http://s3.amazonaws.com/rapgenius/lazy-eye-treatment.jpg", "Success", "None", 41397, "Brian Nguyen"),
    @(37, "http://www.thehealthage.com/site/wp-content/uploads/2010/10/lazy-eye.jpg", "Success", "None", 41397, "Brian Nguyen"),
    @(38, "http://t3.gstatic.com/images?q=tbn:ANd9GcSpfRlvwWWm3T6lSmwpq_a09LVgrB92POV8yjeVhp4awRINl--H", "Error", "Only detect left eye", 41397, "Brian Nguyen")
)

# Give the new date cells (col D) the same date-number-format + wrap style
# already used by the existing D2:D20 cells, by copying that formatting down
# before the values are written.
$ws.Range("D20").Copy() | Out-Null
$ws.Range("D21:D38").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

foreach ($r in $rows) {
    $rowNum = $r[0]
    $url = $r[1]
    $status = $r[2]
    $desc = $r[3]
    $dateSerial = $r[4]
    $updater = $r[5]

    $ws.Cells.Item($rowNum, 1).Value2 = $url
    $ws.Hyperlinks.Add($ws.Cells.Item($rowNum, 1), $url) | Out-Null

    $ws.Cells.Item($rowNum, 2).Value2 = $status
    $ws.Cells.Item($rowNum, 3).Value2 = $desc
    $ws.Cells.Item($rowNum, 4).Value2 = $dateSerial
    $ws.Cells.Item($rowNum, 5).Value2 = $updater
}

# Row heights: Excel auto-fit these to their (now shorter / wrapped) content
# when the rows were edited; row 20's wraps also changed as a side effect.
$ws.Rows.Item(20).RowHeight = 45.75
$ws.Rows.Item(21).RowHeight = 39.75
$ws.Rows.Item(22).RowHeight = 39

# Leave the view where the editor finished: scrolled down with E38 selected.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 22
$ws.Range("E38").Select()
